$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 51f9c755-... and 5b12891f-... rows moves from 04:17:11 to 04:17:59
$wsOverview.Range("G3").Value = "2016-08-16 04:17:59"
$wsOverview.Range("G4").Value = "2016-08-16 04:17:59"

# zh-cn sheet: Priority "ht" -> "mt" for both rows
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime 04:17:00 -> 04:17:53
$wsZhCn.Range("H3").Value = "2016-08-16 04:17:53"
$wsZhCn.Range("H4").Value = "2016-08-16 04:17:53"

# zh-cn sheet: Correspond Handback DateTime 04:17:27 -> 04:18:21
$wsZhCn.Range("K3").Value = "2016-08-16 04:18:21"
$wsZhCn.Range("K4").Value = "2016-08-16 04:18:21"

# de-de sheet: Priority "ht" -> "mt" for both rows
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handoff Datetime 04:17:11 -> 04:17:59
# (shares the same underlying text as Overview's Latest HO Xliff Generate Date)
$wsDeDe.Range("H3").Value = "2016-08-16 04:17:59"
$wsDeDe.Range("H4").Value = "2016-08-16 04:17:59"

# de-de sheet: Correspond Handback DateTime 04:17:35 -> 04:18:28
$wsDeDe.Range("K3").Value = "2016-08-16 04:18:28"
$wsDeDe.Range("K4").Value = "2016-08-16 04:18:28"
